# Apply updated symbol list values (Price / Volume(1h) columns).
# Cells are stored as text in the sheet (t="inlineStr"), so each target
# cell is pre-formatted as Text before the write -- otherwise COM/Excel
# would auto-convert a numeric-looking string like "335.86" or "2.14%"
# into a real number/percentage, which would change the stored type and
# lose exact text formatting (e.g. trailing zeros like "43.90").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '335.86'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '2.14%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '43.90'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '6.26%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.816'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '2.92%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08341'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2.06%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '8.796'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.57%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.984'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-1.51%'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-1.72%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9395'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '1.98%'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.52%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09669'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '3.07%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.04583'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '20.64%'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.84%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001297'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.59%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006031'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-1.52%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.496'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.50%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.505'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.18%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.771'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '5.19%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1362'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.77%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04399'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.23%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004392'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '1.64%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001261'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '4.89%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003991'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02807'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-0.03%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05723'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '5.79%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007922'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '3.12%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1430'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '1.03%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.008967'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '0.14%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002156'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-0.76%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01055'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-8.38%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00007210'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '9.52%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000750'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.09%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003243'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '1.30%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.09%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002101'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.09%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002001'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.09%'
